# HighLevelSequenceDiagrams.pptx edit:
#   "updated HighLevelSequenceDiagram to replace "AddressBook" with "3VIA""
#
# Also renames deletePerson(p) -> deleteCard(c), and (best effort) adds the
# two slide guides recorded in the author's session.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------------
# 1) TextBox 28 (shape 16): "deletePerson(p)" -> "deleteCard(c)"
#    Two runs in the original: "deletePerson" | "(p)". Replace each run's
#    text in place (via Characters over the original run span) so existing
#    run-level formatting is preserved.
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(16)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 12).Text = "deleteCard"   # was "deletePerson"
$tr.Characters(11, 3).Text = "(c)"          # was "(p)"

# ---------------------------------------------------------------------------
# 2) TextBox 32 (shape 17): "post(AddressBookChangedEvent)" ->
#    "post(3VIAChangedEvent)"
#    Original runs: "post(" | "AddressBookChangedEvent" | ")". Collapse the
#    first two runs' span into the new combined text, leaving the trailing
#    ")" run untouched.
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(17)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 28).Text = "post(3VIAChangedEvent"

# ---------------------------------------------------------------------------
# 3) TextBox 61 (shape 29): same replacement as above (second occurrence).
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(29)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 28).Text = "post(3VIAChangedEvent"

# ---------------------------------------------------------------------------
# 4) TextBox 73 (shape 37): "handleAddresssBookChangedEvent()" ->
#    "handle3VIAChangedEvent()"
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(37)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 30).Text = "handle3VIAChangedEvent"

# ---------------------------------------------------------------------------
# 5) TextBox 49 (shape 44): same replacement as above (second occurrence).
# ---------------------------------------------------------------------------
$shp = $s.Shapes.Item(44)
$tr = $shp.TextFrame.TextRange
$tr.Characters(1, 30).Text = "handle3VIAChangedEvent"

# ---------------------------------------------------------------------------
# 6) Presentation-level slide guides (best effort - matches the author's
#    PowerPoint session: one horizontal guide at 186pt, one vertical guide
#    at 360pt).
# ---------------------------------------------------------------------------
try {
    $p.Guides.Add(1, 186)
    $p.Guides.Add(2, 360)
} catch {
    # Guides collection may be unavailable in this environment; ignore.
}
